$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: new columns I (I0) and J (IF) ---
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Match the header formatting (bold, centered, bordered) already used by
# the other header cells (e.g. H1) by copying its format onto I1:J1.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# --- Data rows 2-65: new numeric columns I and J ---
$data = @(
    @(8, 7),
    @(6, 7),
    @(7, 7),
    @(4, 4),
    @(7, 7),
    @(8, 8),
    @(7, 7),
    @(6, 6),
    @(9, 9),
    @(7, 8),
    @(6, 6),
    @(9, 9),
    @(6, 6),
    @(7, 7),
    @(6, 6),
    @(6, 6),
    @(8, 8),
    @(7, 7),
    @(9, 9),
    @(8, 8),
    @(5, 5),
    @(6, 6),
    @(6, 6),
    @(9, 9),
    @(8, 8),
    @(6, 7),
    @(10, 10),
    @(5, 5),
    @(6, 6),
    @(8, 8),
    @(7, 8),
    @(7, 8),
    @(6, 7),
    @(9, 9),
    @(8, 8),
    @(9, 9),
    @(8, 8),
    @(8, 8),
    @(9, 9),
    @(8, 8),
    @(8, 8),
    @(8, 8),
    @(8, 8),
    @(7, 8),
    @(8, 8),
    @(9, 9),
    @(8, 8),
    @(8, 8),
    @(8, 8),
    @(8, 8),
    @(7, 8),
    @(7, 8),
    @(8, 8),
    @(7, 7),
    @(8, 8),
    @(8, 8),
    @(8, 8),
    @(8, 8),
    @(8, 8),
    @(8, 8),
    @(5, 5),
    @(5, 5),
    @(5, 5),
    @(4, 4)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
